# Update the "TFEC" technology list:
#  - keep Residential / Commercial / existing Transport rows as-is
#  - insert a new "Motorcycles" transport row before the Railway rows
#  - explode the single "Industrial uses" row into 5 sub-sector rows
#  - keep the "Other electricity" row as the last row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Residential electricity", "RES_COOL_001", "Residential", "Cooling", "Electricity"),
    @("Residential electricity", "RES_COOL_002", "Residential", "Cooling new users", "Electricity"),
    @("Residential electricity", "RES_CWH_EL_001", "Residential", "Cooking and water heating", "Electricity"),
    @("Residential electricity", "RES_EL_APP_001", "Residential", "Residential appliances", "Electricity"),
    @("Residential electricity", "RES_EL_APP_002", "Residential", "Residential appliances new users", "Electricity"),
    @("Commercial electricity", "COM_EL_APP_001", "Commercial", "Commercial uses", "Electricity"),
    @("Transport electricity", "TRA_BUS_ELC_001", "Transport", "Buses", "Electricity"),
    @("Transport electricity", "TRA_CAR_ELC_001", "Transport", "Cars", "Electricity"),
    @("Transport electricity", "TRA_MCY_EL_001", "Transport", "Motorcycles", "Electricity"),
    @("Transport electricity", "TRA_RLW_FREIGHT_ELC_001", "Transport", "Railway", "Electricity"),
    @("Transport electricity", "TRA_RLW_PSNG_ELC_001", "Transport", "Railway", "Electricity"),
    @("Industrial electricity", "IND_CHEM_ELC_001", "Industry", "Industry chemical", "Electricity"),
    @("Industrial electricity", "IND_CMNT_ELC_001", "Industry", "Industry cement", "Electricity"),
    @("Industrial electricity", "IND_FOOD_ELC_001", "Industry", "Industry food", "Electricity"),
    @("Industrial electricity", "IND_OTH_ELC_001", "Industry", "Industry other", "Electricity"),
    @("Industrial electricity", "IND_TEXT_ELC_001", "Industry", "Industry textile", "Electricity"),
    @("Other electricity", "OTH_ALL_EL_001", "Others", "Other uses", "Electricity")
)

# Template row used to copy the cell formatting (fill/border) for any
# brand-new rows that fall outside the sheet's original used range.
$ws.Range("A2:E2").Copy()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    if ($r -gt 13) {
        $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
    }

    $r = $r + 1
}

[void]$ws.Range("A10").Select()
